# feat: add 2022-Q3 data
#
# The workbook has sheets: 总计, 2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3.
# This script:
#   1. Inserts a new "2022-Q3" sheet (a duplicate of the "2022-Q2" sheet,
#      then rewritten with the new quarter's fund holdings) right after
#      "总计" and before "2022-Q2".
#   2. Inserts a new summary row at the top of the "总计" sheet's data
#      (row 2) describing the 2022-Q3 quarter (13 funds held, 0.79 billion
#      yuan), pushing the existing summary rows down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet by duplicating "2022-Q2"
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")
$wsQ2 = $wb.Worksheets.Item("2022-Q2")

# Places the duplicate immediately before "2022-Q2" (i.e. right after 总计)
$wsQ2.Copy($wsQ2)
$ws3 = $wb.Worksheets.Item(2)
$ws3.Name = "2022-Q3"

# Update the two funds that were already present in 2022-Q2 with their
# refreshed 2022-Q3 figures (fund code/name stay the same).
$ws3.Range("D2:G3").NumberFormat = "@"
$ws3.Range("D2").Value = "11.73"
$ws3.Range("E2").Value = "90.54"
$ws3.Range("F2").Value = "2.88"
$ws3.Range("G2").Value = "0.3378"
$ws3.Range("H2").Value = 7

$ws3.Range("D3").Value = "9.94"
$ws3.Range("E3").Value = "91.00"
$ws3.Range("F3").Value = "2.88"
$ws3.Range("G3").Value = "0.2863"
$ws3.Range("H3").Value = 7

# Append the additional funds newly reported for 2022-Q3.
$ws3.Range("A2").Copy($ws3.Range("A4"))
$ws3.Range("A4").Value = 2
$ws3.Range("B4:G4").NumberFormat = "@"
$ws3.Range("B4").Value = "013603"
$ws3.Range("C4").Value = "易方达均衡优选一年持有混合A"
$ws3.Range("D4").Value = "2.48"
$ws3.Range("E4").Value = "49.36"
$ws3.Range("F4").Value = "1.72"
$ws3.Range("G4").Value = "0.0427"
$ws3.Range("H4").Value = 8

$ws3.Range("A2").Copy($ws3.Range("A5"))
$ws3.Range("A5").Value = 3
$ws3.Range("B5:G5").NumberFormat = "@"
$ws3.Range("B5").Value = "012426"
$ws3.Range("C5").Value = "南方价值臻选混合A"
$ws3.Range("D5").Value = "3.91"
$ws3.Range("E5").Value = "63.50"
$ws3.Range("F5").Value = "0.83"
$ws3.Range("G5").Value = "0.0325"
$ws3.Range("H5").Value = 4

$ws3.Range("A2").Copy($ws3.Range("A6"))
$ws3.Range("A6").Value = 4
$ws3.Range("B6:G6").NumberFormat = "@"
$ws3.Range("B6").Value = "001181"
$ws3.Range("C6").Value = "南方改革机遇灵活配置混合"
$ws3.Range("D6").Value = "3.40"
$ws3.Range("E6").Value = "61.41"
$ws3.Range("F6").Value = "0.82"
$ws3.Range("G6").Value = "0.0279"
$ws3.Range("H6").Value = 4

$ws3.Range("A2").Copy($ws3.Range("A7"))
$ws3.Range("A7").Value = 5
$ws3.Range("B7:G7").NumberFormat = "@"
$ws3.Range("B7").Value = "001536"
$ws3.Range("C7").Value = "南方君选灵活配置混合"
$ws3.Range("D7").Value = "3.38"
$ws3.Range("E7").Value = "56.86"
$ws3.Range("F7").Value = "0.82"
$ws3.Range("G7").Value = "0.0277"
$ws3.Range("H7").Value = 3

$ws3.Range("A2").Copy($ws3.Range("A8"))
$ws3.Range("A8").Value = 6
$ws3.Range("B8:G8").NumberFormat = "@"
$ws3.Range("B8").Value = "202213"
$ws3.Range("C8").Value = "南方核心竞争混合"
$ws3.Range("D8").Value = "2.05"
$ws3.Range("E8").Value = "62.82"
$ws3.Range("F8").Value = "0.82"
$ws3.Range("G8").Value = "0.0168"
$ws3.Range("H8").Value = 6

$ws3.Range("A2").Copy($ws3.Range("A9"))
$ws3.Range("A9").Value = 7
$ws3.Range("B9:G9").NumberFormat = "@"
$ws3.Range("B9").Value = "013604"
$ws3.Range("C9").Value = "易方达均衡优选一年持有混合C"
$ws3.Range("D9").Value = "0.35"
$ws3.Range("E9").Value = "49.36"
$ws3.Range("F9").Value = "1.72"
$ws3.Range("G9").Value = "0.0060"
$ws3.Range("H9").Value = 8

$ws3.Range("A2").Copy($ws3.Range("A10"))
$ws3.Range("A10").Value = 8
$ws3.Range("B10:G10").NumberFormat = "@"
$ws3.Range("B10").Value = "013590"
$ws3.Range("C10").Value = "南方比较优势混合A"
$ws3.Range("D10").Value = "0.46"
$ws3.Range("E10").Value = "66.06"
$ws3.Range("F10").Value = "0.94"
$ws3.Range("G10").Value = "0.0043"
$ws3.Range("H10").Value = 3

$ws3.Range("A2").Copy($ws3.Range("A11"))
$ws3.Range("A11").Value = 9
$ws3.Range("B11:G11").NumberFormat = "@"
$ws3.Range("B11").Value = "013166"
$ws3.Range("C11").Value = "东兴宸祥量化混合A"
$ws3.Range("D11").Value = "0.38"
$ws3.Range("E11").Value = "93.87"
$ws3.Range("F11").Value = "1.11"
$ws3.Range("G11").Value = "0.0042"
$ws3.Range("H11").Value = 10

$ws3.Range("A2").Copy($ws3.Range("A12"))
$ws3.Range("A12").Value = 10
$ws3.Range("B12:G12").NumberFormat = "@"
$ws3.Range("B12").Value = "013591"
$ws3.Range("C12").Value = "南方比较优势混合C"
$ws3.Range("D12").Value = "0.30"
$ws3.Range("E12").Value = "66.06"
$ws3.Range("F12").Value = "0.94"
$ws3.Range("G12").Value = "0.0028"
$ws3.Range("H12").Value = 3

$ws3.Range("A2").Copy($ws3.Range("A13"))
$ws3.Range("A13").Value = 11
$ws3.Range("B13:G13").NumberFormat = "@"
$ws3.Range("B13").Value = "012427"
$ws3.Range("C13").Value = "南方价值臻选混合C"
$ws3.Range("D13").Value = "0.19"
$ws3.Range("E13").Value = "63.50"
$ws3.Range("F13").Value = "0.83"
$ws3.Range("G13").Value = "0.0016"
$ws3.Range("H13").Value = 4

$ws3.Range("A2").Copy($ws3.Range("A14"))
$ws3.Range("A14").Value = 12
$ws3.Range("B14:G14").NumberFormat = "@"
$ws3.Range("B14").Value = "013167"
$ws3.Range("C14").Value = "东兴宸祥量化混合C"
$ws3.Range("D14").Value = "0.08"
$ws3.Range("E14").Value = "93.87"
$ws3.Range("F14").Value = "1.11"
$ws3.Range("G14").Value = "0.0009"
$ws3.Range("H14").Value = 10

# ---------------------------------------------------------------------
# Step 2: insert the 2022-Q3 summary row into "总计"
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()
$wsTotal.Range("A3").Copy($wsTotal.Range("A2"))
$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 13
$wsTotal.Range("D2").Value = 0.79

# The existing rows kept their old running index when they shifted down a
# row, so renumber the sequential index column to stay 0,1,2,3,4.
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("A4").Value = 2
$wsTotal.Range("A5").Value = 3
$wsTotal.Range("A6").Value = 4

# Keep the originally-active tab ("2021-Q3") selected, since duplicating
# the "2022-Q2" sheet above made the new sheet the active one.
$wb.Worksheets.Item("2021-Q3").Activate()
